$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3623.75
$ws.Range("I40").Value = 3747.5
$ws.Range("K40").Value = 3747.5
$ws.Range("M40").Value = -3572.5

# row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 3780.15
$ws.Range("I138").Value = 4499.5
$ws.Range("J138").Value = 3700.2222
$ws.Range("K138").Value = 13498.5
$ws.Range("L138").Value = 11100.6666
$ws.Range("M138").Value = -8358.5
$ws.Range("N138").Value = -21380.6666

$ws = $wb.Worksheets.Item("ARM")
# row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 31000
$ws.Range("I37").Value = 24000
$ws.Range("K37").Value = 24000
$ws.Range("M37").Value = -23727

# row 55 (Leve Item ID 2830)
$ws.Range("H55").Value = 48500
$ws.Range("J55").Value = 48500
$ws.Range("L55").Value = 48500
$ws.Range("N55").Value = -49130

# row 92 (Leve Item ID 18050)
$ws.Range("H92").Value = 41912.5
$ws.Range("J92").Value = 41912.5
$ws.Range("L92").Value = 41912.5
$ws.Range("N92").Value = -46904.5

# row 94 (Leve Item ID 18055)
$ws.Range("H94").Value = 59999
$ws.Range("J94").Value = 59999
$ws.Range("L94").Value = 59999
$ws.Range("N94").Value = -61801

# row 107 (Leve Item ID 25645)
$ws.Range("H107").Value = 25000
$ws.Range("J107").Value = 25000
$ws.Range("L107").Value = 25000
$ws.Range("N107").Value = -32680

# row 109 (Leve Item ID 25646)
$ws.Range("H109").Value = 99500
$ws.Range("J109").Value = 99500
$ws.Range("L109").Value = 99500
$ws.Range("N109").Value = -102274

# row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1609.7273
$ws.Range("I132").Value = 1378.5555
$ws.Range("K132").Value = 4135.666499999999
$ws.Range("M132").Value = -1605.666499999999

$ws = $wb.Worksheets.Item("BSM")
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1979.3636
$ws.Range("I86").Value = 1863.7778
$ws.Range("K86").Value = 1863.7778
$ws.Range("M86").Value = -740.7778000000001

# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1979.3636
$ws.Range("I89").Value = 1863.7778
$ws.Range("K89").Value = 9318.889000000001
$ws.Range("M89").Value = -3702.889000000001

$ws = $wb.Worksheets.Item("CRP")
# row 53 (Leve Item ID 25632)
$ws.Range("H53").Value = 30999
$ws.Range("J53").Value = 30999
$ws.Range("L53").Value = 30999
$ws.Range("N53").Value = -32213

# row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 25998.334
$ws.Range("I60").Value = 14995
$ws.Range("J60").Value = 31500
$ws.Range("K60").Value = 14995
$ws.Range("L60").Value = 31500
$ws.Range("M60").Value = -14484
$ws.Range("N60").Value = -32522

# row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 2199.5
$ws.Range("I62").Value = 2199
$ws.Range("J62").Value = 2200
$ws.Range("K62").Value = 2199
$ws.Range("L62").Value = 2200
$ws.Range("M62").Value = -1575
$ws.Range("N62").Value = -3448

# row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 2199.5
$ws.Range("I65").Value = 2199
$ws.Range("J65").Value = 2200
$ws.Range("K65").Value = 10995
$ws.Range("L65").Value = 11000
$ws.Range("M65").Value = -7875
$ws.Range("N65").Value = -17240

# row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 34000
$ws.Range("I68").Value = 34000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 34000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -33251
$ws.Range("N68").ClearContents()

# row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 34000
$ws.Range("I71").Value = 34000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 102000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -98256
$ws.Range("N71").ClearContents()

# row 74 (Leve Item ID 10636)
$ws.Range("H74").Value = 15314
$ws.Range("J74").Value = 15314
$ws.Range("L74").Value = 15314
$ws.Range("N74").Value = -17062

# row 77 (Leve Item ID 10636)
$ws.Range("H77").Value = 15314
$ws.Range("J77").Value = 15314
$ws.Range("L77").Value = 45942
$ws.Range("N77").Value = -54678

# row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 13945594
$ws.Range("I86").Value = 17430742
$ws.Range("K86").Value = 17430742
$ws.Range("M86").Value = -17429619

# row 88 (Leve Item ID 10608)
$ws.Range("H88").Value = 21171.5
$ws.Range("J88").Value = 21171.5
$ws.Range("L88").Value = 21171.5
$ws.Range("N88").Value = -21983.5

# row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 13945594
$ws.Range("I89").Value = 17430742
$ws.Range("K89").Value = 87153710
$ws.Range("M89").Value = -87148094

# row 91 (Leve Item ID 10608)
$ws.Range("H91").Value = 21171.5
$ws.Range("J91").Value = 21171.5
$ws.Range("L91").Value = 21171.5
$ws.Range("N91").Value = -23979.5

$ws = $wb.Worksheets.Item("CUL")
# row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 986.93335
$ws.Range("I55").Value = 776
$ws.Range("J55").Value = 1019.38464
$ws.Range("K55").Value = 2328
$ws.Range("L55").Value = 3058.15392
$ws.Range("M55").Value = -2151
$ws.Range("N55").Value = -3412.15392

$ws = $wb.Worksheets.Item("GSM")
# row 43 (Leve Item ID 4218)
$ws.Range("H43").Value = 3282.2856
$ws.Range("J43").Value = 3999.5
$ws.Range("L43").Value = 3999.5
$ws.Range("N43").Value = -4301.5

# row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 11090.857
$ws.Range("J46").Value = 16059
$ws.Range("L46").Value = 16059
$ws.Range("N46").Value = -16371

# row 55 (Leve Item ID 4237)
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

# row 92 (Leve Item ID 18094)
$ws.Range("H92").Value = 12175.571
$ws.Range("J92").Value = 12175.571
$ws.Range("L92").Value = 12175.571
$ws.Range("N92").Value = -15919.571

# row 94 (Leve Item ID 19511)
$ws.Range("H94").Value = 29999
$ws.Range("J94").Value = 29999
$ws.Range("L94").Value = 29999
$ws.Range("N94").Value = -31351

# row 96 (Leve Item ID 18261)
$ws.Range("H96").Value = 25000
$ws.Range("J96").Value = 25000
$ws.Range("L96").Value = 25000
$ws.Range("N96").Value = -30492

# row 98 (Leve Item ID 18359)
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2674.625
$ws.Range("I100").Value = 2233
$ws.Range("K100").Value = 2233
$ws.Range("M100").Value = -1692

# row 140 (Leve Item ID 42503)
$ws.Range("H140").Value = 79000
$ws.Range("J140").Value = 79000
$ws.Range("L140").Value = 79000
$ws.Range("N140").Value = -89360

$ws = $wb.Worksheets.Item("WVR")
# row 108 (Leve Item ID 25661)
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
